$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new dSF (column F) value, per the target diff
$updates = @{
    3  = 0
    8  = 12
    10 = -4
    11 = -4
    12 = 3
    13 = 5
    14 = -5
    15 = 1
    16 = 5
    18 = -3
    19 = -1
    21 = -2
    22 = -1
    23 = 9
    25 = -2
    26 = -2
    27 = -5
    30 = 1
    31 = -2
    32 = 3
    33 = -3
    34 = -1
    36 = -2
    37 = -6
    38 = 0
    40 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
